$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: new peer-review entry (B. Allton) ---
$ws.Range("A4").Value = "B. Allton"
$ws.Range("B4").Value = "2/14/2020"
$ws.Range("B4").NumberFormat = "d-mmm-yy"
$ws.Range("C4").Value = "Document for grammar and flow "
$ws.Range("D4").Value = "Some sentences in part 7 not quite clear "
$ws.Range("E4").Value = "Reworded some of the paragraphs to give them better flow and more clarification as to what was done throughout the project history.  "

# Row 4 alignment: A/B left+center, C center+center+wrap, D/E wrap only
$ws.Range("A4:B4").HorizontalAlignment = -4131
$ws.Range("A4:B4").VerticalAlignment = -4108

$ws.Range("C4").HorizontalAlignment = -4108
$ws.Range("C4").VerticalAlignment = -4108
$ws.Range("C4").WrapText = $true

$ws.Range("D4:E4").WrapText = $true

# Row 4 custom height to fit the wrapped text
$ws.Range("A4:E4").RowHeight = 90

# --- Rows 5-15: apply left/center alignment to A/B, wrap to C/D/E ---
$ws.Range("A5:B15").HorizontalAlignment = -4131
$ws.Range("A5:B15").VerticalAlignment = -4108

$ws.Range("C5:E15").WrapText = $true

# --- Restore the last on-screen selection ---
$ws.Range("K4:K7").Select()
